$wb = $excel.ActiveWorkbook

# Template sheet to copy from (Norway already has the expected layout/styles)
$source = $wb.Worksheets.Item("Norway")

# New country sheets, in tab order, with the B2 (market name) and B4 (ticket id) values.
# TicketFirst controls whether B4 is written before B2 (matches the order the
# shared strings were originally authored in, so the shared-string table
# indices line up with the source diff).
$countries = @(
    @{ Name = "Hungary";  Market = "Hungary Market";  Ticket = $null;              Active = "A1:XFD1048576"; TicketFirst = $false },
    @{ Name = "Italy";    Market = "Italy market";     Ticket = "NGC-3443/T1916";   Active = "B6";            TicketFirst = $false },
    @{ Name = "Spain";    Market = "Spain market";     Ticket = "NGC-3442/T1592";   Active = "A1:XFD1048576"; TicketFirst = $false },
    @{ Name = "Serbia";   Market = "Serbia market";    Ticket = "NGC-4305/T3495";   Active = "B6";            TicketFirst = $false },
    @{ Name = "Romania";  Market = "Romania market";   Ticket = "NGC-4307/T3541";   Active = "A1:XFD1048576"; TicketFirst = $false },
    @{ Name = "Slovakia"; Market = "Slovakia market";  Ticket = "NGC-4306/T3555";   Active = "B4";            TicketFirst = $true },
    @{ Name = "Turkey";   Market = "Turkey Market";    Ticket = $null;              Active = "B5";            TicketFirst = $false }
)

$last = $source
foreach ($c in $countries) {
    $source.Copy($null, $last)
    $newSheet = $wb.Worksheets.Item($last.Index + 1)
    $newSheet.Name = $c.Name

    if ($c.TicketFirst) {
        if ($c.Ticket -ne $null) {
            $newSheet.Range("B4").Value = $c.Ticket
        } else {
            $newSheet.Range("B4").Value = ""
        }
        $newSheet.Range("B2").Value = $c.Market
    } else {
        $newSheet.Range("B2").Value = $c.Market
        if ($c.Ticket -ne $null) {
            $newSheet.Range("B4").Value = $c.Ticket
        } else {
            $newSheet.Range("B4").Value = ""
        }
    }

    $newSheet.Range($c.Active).Select() | Out-Null

    $last = $newSheet
}

$last.Select() | Out-Null
